$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update team-specific time transition matrix values (row by row)
# Row 2
$ws.Range("B2").Value = 0.2258064516129032
$ws.Range("C2").Value = 0.5080645161290323
$ws.Range("J2").Value = 0.004032258064516129
$ws.Range("O2").Value = 0.004032258064516129
$ws.Range("P2").Value = 0.157258064516129
$ws.Range("S2").Value = 0.1008064516129032

# Row 3
$ws.Range("C3").Value = 0.03053435114503817
$ws.Range("P3").Value = 0.7557251908396947
$ws.Range("S3").Value = 0.2137404580152672

# Row 4
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.7058823529411765
$ws.Range("S4").Value = 0.2352941176470588

# Row 6
$ws.Range("B6").Value = 0.07692307692307693
$ws.Range("D6").Value = 0.005128205128205128
$ws.Range("F6").Value = 0.04615384615384616
$ws.Range("J6").Value = 0.3076923076923077
$ws.Range("O6").Value = 0.02564102564102564
$ws.Range("Q6").Value = 0.1743589743589744
$ws.Range("R6").Value = 0.1333333333333333
$ws.Range("S6").Value = 0.2307692307692308

# Row 7
$ws.Range("B7").Value = 0.06790123456790123
$ws.Range("D7").Value = 0.006172839506172839
$ws.Range("F7").Value = 0.03703703703703703
$ws.Range("J7").Value = 0.09876543209876543
$ws.Range("O7").Value = 0.01851851851851852
$ws.Range("Q7").Value = 0.1975308641975309
$ws.Range("R7").Value = 0.1172839506172839
$ws.Range("S7").Value = 0.4567901234567901

# Row 8
$ws.Range("B8").Value = 0.06483790523690773
$ws.Range("D8").Value = 0.02493765586034913
$ws.Range("F8").Value = 0.05735660847880299
$ws.Range("J8").Value = 0.1172069825436409
$ws.Range("O8").Value = 0.009975062344139651
$ws.Range("Q8").Value = 0.2319201995012469
$ws.Range("R8").Value = 0.114713216957606
$ws.Range("S8").Value = 0.3790523690773067

# Row 9
$ws.Range("B9").Value = 0.0963302752293578
$ws.Range("D9").Value = 0.03211009174311927
$ws.Range("F9").Value = 0.05963302752293578
$ws.Range("J9").Value = 0.08256880733944955
$ws.Range("O9").Value = 0.02293577981651376
$ws.Range("Q9").Value = 0.2201834862385321
$ws.Range("R9").Value = 0.09174311926605505
$ws.Range("S9").Value = 0.3944954128440367

# Row 10
$ws.Range("B10").Value = 0.09055425448868072
$ws.Range("D10").Value = 0.01327088212334114
$ws.Range("F10").Value = 0.05776736924277908
$ws.Range("J10").Value = 0.1030444964871194
$ws.Range("O10").Value = 0.01249024199843872
$ws.Range("Q10").Value = 0.2318501170960187
$ws.Range("R10").Value = 0.1053864168618267
$ws.Range("S10").Value = 0.3856362217017955

# Row 11
$ws.Range("G11").Value = 0.1341991341991342
$ws.Range("J11").Value = 0.0735930735930736
$ws.Range("K11").Value = 0.1861471861471861
$ws.Range("L11").Value = 0.5974025974025974
$ws.Range("S11").Value = 0.008658008658008658

# Row 12
$ws.Range("G12").Value = 0.8041958041958042
$ws.Range("J12").Value = 0.1468531468531468
$ws.Range("K12").Value = 0.006993006993006993
$ws.Range("L12").Value = 0.04195804195804196

# Row 13
$ws.Range("G13").Value = 0.6896551724137931
$ws.Range("J13").Value = 0.3103448275862069

# Row 15
$ws.Range("F15").Value = 0.04017857142857143
$ws.Range("H15").Value = 0.1205357142857143
$ws.Range("I15").Value = 0.0625
$ws.Range("J15").Value = 0.40625
$ws.Range("K15").Value = 0.03571428571428571
$ws.Range("M15").Value = 0.008928571428571428
$ws.Range("O15").Value = 0.06696428571428571
$ws.Range("S15").Value = 0.2589285714285715

# Row 16
$ws.Range("F16").Value = 0.006369426751592357
$ws.Range("H16").Value = 0.1401273885350318
$ws.Range("I16").Value = 0.1082802547770701
$ws.Range("J16").Value = 0.4840764331210191
$ws.Range("K16").Value = 0.06369426751592357
$ws.Range("M16").Value = 0.01273885350318471
$ws.Range("O16").Value = 0.09554140127388536
$ws.Range("S16").Value = 0.08917197452229299

# Row 17
$ws.Range("F17").Value = 0.01587301587301587
$ws.Range("H17").Value = 0.1825396825396825
$ws.Range("I17").Value = 0.09523809523809523
$ws.Range("J17").Value = 0.4325396825396826
$ws.Range("K17").Value = 0.06944444444444445
$ws.Range("M17").Value = 0.01587301587301587
$ws.Range("O17").Value = 0.07738095238095238
$ws.Range("S17").Value = 0.1111111111111111

# Row 18
$ws.Range("F18").Value = 0.01229508196721311
$ws.Range("H18").Value = 0.1270491803278689
$ws.Range("I18").Value = 0.139344262295082
$ws.Range("J18").Value = 0.4098360655737705
$ws.Range("K18").Value = 0.0860655737704918
$ws.Range("M18").Value = 0.004098360655737705
$ws.Range("O18").Value = 0.08196721311475409
$ws.Range("S18").Value = 0.139344262295082

# Row 19
$ws.Range("F19").Value = 0.01848739495798319
$ws.Range("H19").Value = 0.1907563025210084
$ws.Range("I19").Value = 0.08823529411764706
$ws.Range("J19").Value = 0.4117647058823529
$ws.Range("K19").Value = 0.09411764705882353
$ws.Range("M19").Value = 0.01428571428571429
$ws.Range("N19").Value = 0.001680672268907563
$ws.Range("O19").Value = 0.06386554621848739
$ws.Range("S19").Value = 0.1168067226890756
